# Update Formal Verification Plan - clean-up
# Remove verification items not implemented in cv32e40p:
# mark Coverage Method as "Partial Proof" and add a note in Coverage Details
# explaining what was left open, for all RV32M instruction rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RV32M")

for ($row = 2; $row -le 9; $row++) {
    $ws.Cells.Item($row, 9).Value = "Partial Proof"
    $ws.Cells.Item($row, 10).Value = "details on what was left open"
}

# Update the view state to reflect where the user ended up after editing
$ws.Activate()
$ws.Application.ActiveWindow.Zoom = 66
$null = $ws.Range("J4").Select()
